# Auto-generated: applies the row-content permutation described in the commit diff.
# For each listed row, overwrite columns B and E..AD with the values that, in the
# source data refresh, ended up belonging to that row (columns A, C, D are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("South Africa Premier")

# Row 12
$ws.Range("B12").Value = 7007757
$ws.Range("E12").Value = "Mamelodi Sundowns"
$ws.Range("F12").Value = "Kaizer Chiefs"
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = "H"
$ws.Range("L12").Value = 1.65
$ws.Range("M12").Value = 3.5
$ws.Range("N12").Value = 5
$ws.Range("O12").Value = 1.45
$ws.Range("P12").Value = 4.2
$ws.Range("Q12").Value = 6.5
$ws.Range("R12").Value = -1
$ws.Range("S12").Value = 1.75
$ws.Range("T12").Value = 2.05
$ws.Range("U12").Value = 2.25
$ws.Range("V12").Value = 1.8
$ws.Range("W12").Value = 2
$ws.Range("X12").Value = 0.45
$ws.Range("Y12").Value = -1
$ws.Range("Z12").Value = -1
$ws.Range("AA12").Value = 0
$ws.Range("AB12").Value = 0
$ws.Range("AC12").Value = 0.8
$ws.Range("AD12").Value = -1

# Row 13
$ws.Range("B13").Value = 7007756
$ws.Range("E13").Value = "Moroka Swallows"
$ws.Range("F13").Value = "Cape Town City"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = "A"
$ws.Range("L13").Value = 2.625
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 2.6
$ws.Range("O13").Value = 2.625
$ws.Range("P13").Value = 2.9
$ws.Range("Q13").Value = 2.625
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 1.925
$ws.Range("T13").Value = 1.875
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = 1.975
$ws.Range("W13").Value = 1.825
$ws.Range("X13").Value = -1
$ws.Range("Y13").Value = -1
$ws.Range("Z13").Value = 1.625
$ws.Range("AA13").Value = -1
$ws.Range("AB13").Value = 0.875
$ws.Range("AC13").Value = -1
$ws.Range("AD13").Value = 0.825

# Row 15
$ws.Range("B15").Value = 7007759
$ws.Range("E15").Value = "Polokwane City"
$ws.Range("F15").Value = "Stellenbosch FC"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = "H"
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 3.1
$ws.Range("N15").Value = 2.25
$ws.Range("O15").Value = 2.875
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 2.45
$ws.Range("R15").Value = 0.25
$ws.Range("S15").Value = 1.725
$ws.Range("T15").Value = 2.075
$ws.Range("U15").Value = 2
$ws.Range("V15").Value = 1.85
$ws.Range("W15").Value = 1.95
$ws.Range("X15").Value = 1.875
$ws.Range("Y15").Value = -1
$ws.Range("Z15").Value = -1
$ws.Range("AA15").Value = 0.7250000000000001
$ws.Range("AB15").Value = -1
$ws.Range("AC15").Value = -1
$ws.Range("AD15").Value = 0.95

# Row 16
$ws.Range("B16").Value = 7007760
$ws.Range("E16").Value = "Chippa United"
$ws.Range("F16").Value = "TS Galaxy"
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = "D"
$ws.Range("L16").Value = 2.4
$ws.Range("M16").Value = 2.9
$ws.Range("N16").Value = 3
$ws.Range("O16").Value = 2.4
$ws.Range("P16").Value = 2.8
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = -0.25
$ws.Range("S16").Value = 2.075
$ws.Range("T16").Value = 1.725
$ws.Range("U16").Value = 1.75
$ws.Range("V16").Value = 1.825
$ws.Range("W16").Value = 1.975
$ws.Range("X16").Value = -1
$ws.Range("Y16").Value = 1.8
$ws.Range("Z16").Value = -1
$ws.Range("AA16").Value = -0.5
$ws.Range("AB16").Value = 0.3625
$ws.Range("AC16").Value = -1
$ws.Range("AD16").Value = 0.9750000000000001

# Row 17
$ws.Range("B17").Value = 7007761
$ws.Range("E17").Value = "Ajax Cape Town"
$ws.Range("F17").Value = "Sekhukhune United FC"
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = "A"
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 3
$ws.Range("N17").Value = 2.3
$ws.Range("O17").Value = 2.7
$ws.Range("P17").Value = 2.625
$ws.Range("Q17").Value = 2.75
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 1.9
$ws.Range("T17").Value = 1.9
$ws.Range("U17").Value = 2
$ws.Range("V17").Value = 2.025
$ws.Range("W17").Value = 1.775
$ws.Range("X17").Value = -1
$ws.Range("Y17").Value = -1
$ws.Range("Z17").Value = 1.75
$ws.Range("AA17").Value = -1
$ws.Range("AB17").Value = 0.8999999999999999
$ws.Range("AC17").Value = 0
$ws.Range("AD17").Value = 0

# Row 64
$ws.Range("B64").Value = 7267513
$ws.Range("E64").Value = "Richards Bay FC"
$ws.Range("F64").Value = "Polokwane City"
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 1
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1
$ws.Range("K64").Value = "A"
$ws.Range("L64").Value = 2.55
$ws.Range("M64").Value = 2.875
$ws.Range("N64").Value = 2.875
$ws.Range("O64").Value = 2.2
$ws.Range("P64").Value = 2.875
$ws.Range("Q64").Value = 3.6
$ws.Range("R64").Value = -0.25
$ws.Range("S64").Value = 1.925
$ws.Range("T64").Value = 1.875
$ws.Range("U64").Value = 2
$ws.Range("V64").Value = 2.025
$ws.Range("W64").Value = 1.775
$ws.Range("X64").Value = -1
$ws.Range("Y64").Value = -1
$ws.Range("Z64").Value = 2.6
$ws.Range("AA64").Value = -1
$ws.Range("AB64").Value = 0.875
$ws.Range("AC64").Value = -1
$ws.Range("AD64").Value = 0.7749999999999999

# Row 65
$ws.Range("B65").Value = 7267510
$ws.Range("E65").Value = "Stellenbosch FC"
$ws.Range("F65").Value = "TS Galaxy"
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 2
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = "H"
$ws.Range("L65").Value = 2.35
$ws.Range("M65").Value = 2.875
$ws.Range("N65").Value = 3.2
$ws.Range("O65").Value = 2.3
$ws.Range("P65").Value = 2.875
$ws.Range("Q65").Value = 3.3
$ws.Range("R65").Value = -0.25
$ws.Range("S65").Value = 1.975
$ws.Range("T65").Value = 1.825
$ws.Range("U65").Value = 2
$ws.Range("V65").Value = 1.85
$ws.Range("W65").Value = 1.95
$ws.Range("X65").Value = 1.3
$ws.Range("Y65").Value = -1
$ws.Range("Z65").Value = -1
$ws.Range("AA65").Value = 0.9750000000000001
$ws.Range("AB65").Value = -1
$ws.Range("AC65").Value = 0.8500000000000001
$ws.Range("AD65").Value = -1

# Row 92
$ws.Range("B92").Value = 7404450
$ws.Range("E92").Value = "TS Galaxy"
$ws.Range("F92").Value = "Polokwane City"
$ws.Range("G92").Value = 3
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = "H"
$ws.Range("L92").Value = 3
$ws.Range("M92").Value = 2.875
$ws.Range("N92").Value = 2.4
$ws.Range("O92").Value = 2.5
$ws.Range("P92").Value = 2.75
$ws.Range("Q92").Value = 2.9
$ws.Range("R92").Value = 0
$ws.Range("S92").Value = 1.7
$ws.Range("T92").Value = 2.1
$ws.Range("U92").Value = 1.75
$ws.Range("V92").Value = 1.75
$ws.Range("W92").Value = 2.05
$ws.Range("X92").Value = 1.5
$ws.Range("Y92").Value = -1
$ws.Range("Z92").Value = -1
$ws.Range("AA92").Value = 0.7
$ws.Range("AB92").Value = -1
$ws.Range("AC92").Value = 0.75
$ws.Range("AD92").Value = -1

# Row 93
$ws.Range("B93").Value = 7404449
$ws.Range("E93").Value = "Moroka Swallows"
$ws.Range("F93").Value = "Kaizer Chiefs"
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 1
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = "A"
$ws.Range("L93").Value = 2.45
$ws.Range("M93").Value = 2.9
$ws.Range("N93").Value = 2.9
$ws.Range("O93").Value = 2.75
$ws.Range("P93").Value = 2.7
$ws.Range("Q93").Value = 2.8
$ws.Range("R93").Value = 0
$ws.Range("S93").Value = 1.85
$ws.Range("T93").Value = 1.95
$ws.Range("U93").Value = 2
$ws.Range("V93").Value = 2.05
$ws.Range("W93").Value = 1.75
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z93").Value = 1.8
$ws.Range("AA93").Value = -1
$ws.Range("AB93").Value = 0.95
$ws.Range("AC93").Value = -1
$ws.Range("AD93").Value = 0.75

# Row 105
$ws.Range("B105").Value = 7464395
$ws.Range("E105").Value = "Mamelodi Sundowns"
$ws.Range("F105").Value = "Ajax Cape Town"
$ws.Range("G105").Value = 3
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 2
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = "H"
$ws.Range("L105").Value = 1.181
$ws.Range("M105").Value = 6.5
$ws.Range("N105").Value = 13
$ws.Range("O105").Value = 1.2
$ws.Range("P105").Value = 6
$ws.Range("Q105").Value = 12
$ws.Range("R105").Value = -1.75
$ws.Range("S105").Value = 1.775
$ws.Range("T105").Value = 2.025
$ws.Range("U105").Value = 3
$ws.Range("V105").Value = 1.975
$ws.Range("W105").Value = 1.825
$ws.Range("X105").Value = 0.2
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = 0.7749999999999999
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 0
$ws.Range("AD105").Value = 0

# Row 106
$ws.Range("B106").Value = 7476568
$ws.Range("E106").Value = "Amazulu"
$ws.Range("F106").Value = "Sekhukhune United FC"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 1
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = "H"
$ws.Range("L106").Value = 2.2
$ws.Range("M106").Value = 3
$ws.Range("N106").Value = 3.4
$ws.Range("O106").Value = 2.25
$ws.Range("P106").Value = 2.75
$ws.Range("Q106").Value = 3.6
$ws.Range("R106").Value = -0.25
$ws.Range("S106").Value = 1.9
$ws.Range("T106").Value = 1.9
$ws.Range("U106").Value = 1.75
$ws.Range("V106").Value = 1.975
$ws.Range("W106").Value = 1.825
$ws.Range("X106").Value = 1.25
$ws.Range("Y106").Value = -1
$ws.Range("Z106").Value = -1
$ws.Range("AA106").Value = 0.8999999999999999
$ws.Range("AB106").Value = -1
$ws.Range("AC106").Value = -1
$ws.Range("AD106").Value = 0.825

# Row 111
$ws.Range("B111").Value = 7477991
$ws.Range("E111").Value = "Royal AM FC"
$ws.Range("F111").Value = "Moroka Swallows"
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = "D"
$ws.Range("L111").Value = 2.75
$ws.Range("M111").Value = 2.9
$ws.Range("N111").Value = 2.625
$ws.Range("O111").Value = 2.625
$ws.Range("P111").Value = 2.9
$ws.Range("Q111").Value = 2.7
$ws.Range("R111").Value = 0
$ws.Range("S111").Value = 1.9
$ws.Range("T111").Value = 1.9
$ws.Range("U111").Value = 2
$ws.Range("V111").Value = 1.975
$ws.Range("W111").Value = 1.825
$ws.Range("X111").Value = -1
$ws.Range("Y111").Value = 1.9
$ws.Range("Z111").Value = -1
$ws.Range("AA111").Value = 0
$ws.Range("AB111").Value = 0
$ws.Range("AC111").Value = -1
$ws.Range("AD111").Value = 0.825

# Row 112
$ws.Range("B112").Value = 7477992
$ws.Range("E112").Value = "Chippa United"
$ws.Range("F112").Value = "Polokwane City"
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 2
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2
$ws.Range("K112").Value = "A"
$ws.Range("L112").Value = 2.375
$ws.Range("M112").Value = 3
$ws.Range("N112").Value = 3.1
$ws.Range("O112").Value = 2.5
$ws.Range("P112").Value = 2.9
$ws.Range("Q112").Value = 3
$ws.Range("R112").Value = 0
$ws.Range("S112").Value = 1.775
$ws.Range("T112").Value = 2.025
$ws.Range("U112").Value = 1.75
$ws.Range("V112").Value = 1.75
$ws.Range("W112").Value = 2.05
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = -1
$ws.Range("Z112").Value = 2
$ws.Range("AA112").Value = -1
$ws.Range("AB112").Value = 1.025
$ws.Range("AC112").Value = 0.75
$ws.Range("AD112").Value = -1

# Row 113
$ws.Range("B113").Value = 7478004
$ws.Range("E113").Value = "Kaizer Chiefs"
$ws.Range("F113").Value = "Richards Bay FC"
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 1
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = "H"
$ws.Range("L113").Value = 1.8
$ws.Range("M113").Value = 3
$ws.Range("N113").Value = 5
$ws.Range("O113").Value = 1.85
$ws.Range("P113").Value = 3
$ws.Range("Q113").Value = 4.75
$ws.Range("R113").Value = -0.5
$ws.Range("S113").Value = 1.875
$ws.Range("T113").Value = 1.925
$ws.Range("U113").Value = 1.75
$ws.Range("V113").Value = 1.75
$ws.Range("W113").Value = 2.05
$ws.Range("X113").Value = 0.8500000000000001
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = -1
$ws.Range("AA113").Value = 0.875
$ws.Range("AB113").Value = -1
$ws.Range("AC113").Value = -1
$ws.Range("AD113").Value = 1.05

# Row 118
$ws.Range("B118").Value = 7477996
$ws.Range("E118").Value = "Richards Bay FC"
$ws.Range("F118").Value = "Supersport United"
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 1
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = 1
$ws.Range("K118").Value = "H"
$ws.Range("L118").Value = 4
$ws.Range("M118").Value = 3
$ws.Range("N118").Value = 1.95
$ws.Range("O118").Value = 3.6
$ws.Range("P118").Value = 2.9
$ws.Range("Q118").Value = 2.1
$ws.Range("R118").Value = 0.25
$ws.Range("S118").Value = 1.975
$ws.Range("T118").Value = 1.825
$ws.Range("U118").Value = 2
$ws.Range("V118").Value = 2
$ws.Range("W118").Value = 1.8
$ws.Range("X118").Value = 2.6
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = 0.9750000000000001
$ws.Range("AB118").Value = -1
$ws.Range("AC118").Value = 1
$ws.Range("AD118").Value = -1

# Row 119
$ws.Range("B119").Value = 7477995
$ws.Range("E119").Value = "TS Galaxy"
$ws.Range("F119").Value = "Chippa United"
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = "H"
$ws.Range("L119").Value = 2.3
$ws.Range("M119").Value = 3
$ws.Range("N119").Value = 3.1
$ws.Range("O119").Value = 2.1
$ws.Range("P119").Value = 3
$ws.Range("Q119").Value = 3.6
$ws.Range("R119").Value = -0.25
$ws.Range("S119").Value = 1.825
$ws.Range("T119").Value = 1.975
$ws.Range("U119").Value = 2
$ws.Range("V119").Value = 2.05
$ws.Range("W119").Value = 1.75
$ws.Range("X119").Value = 1.1
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = 0.825
$ws.Range("AB119").Value = -1
$ws.Range("AC119").Value = 0
$ws.Range("AD119").Value = 0

# Row 142
$ws.Range("B142").Value = 7628939
$ws.Range("E142").Value = "Orlando Pirates"
$ws.Range("F142").Value = "Ajax Cape Town"
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 1
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = 0
$ws.Range("K142").Value = "D"
$ws.Range("L142").Value = 1.363
$ws.Range("M142").Value = 4.5
$ws.Range("N142").Value = 7.5
$ws.Range("O142").Value = 1.3
$ws.Range("P142").Value = 5
$ws.Range("Q142").Value = 7.5
$ws.Range("R142").Value = -1.5
$ws.Range("S142").Value = 1.925
$ws.Range("T142").Value = 1.875
$ws.Range("U142").Value = 2.75
$ws.Range("V142").Value = 1.775
$ws.Range("W142").Value = 2.025
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = 4
$ws.Range("Z142").Value = -1
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.875
$ws.Range("AC142").Value = -1
$ws.Range("AD142").Value = 1.025

# Row 143
$ws.Range("B143").Value = 7628938
$ws.Range("E143").Value = "Chippa United"
$ws.Range("F143").Value = "Richards Bay FC"
$ws.Range("G143").Value = 3
$ws.Range("H143").Value = 0
$ws.Range("I143").Value = 2
$ws.Range("J143").Value = 0
$ws.Range("K143").Value = "H"
$ws.Range("L143").Value = 2.375
$ws.Range("M143").Value = 2.75
$ws.Range("N143").Value = 3.1
$ws.Range("O143").Value = 2.25
$ws.Range("P143").Value = 2.75
$ws.Range("Q143").Value = 3.4
$ws.Range("R143").Value = -0.25
$ws.Range("S143").Value = 1.975
$ws.Range("T143").Value = 1.825
$ws.Range("U143").Value = 2
$ws.Range("V143").Value = 2.025
$ws.Range("W143").Value = 1.775
$ws.Range("X143").Value = 1.25
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 0.9750000000000001
$ws.Range("AB143").Value = -1
$ws.Range("AC143").Value = 1.025
$ws.Range("AD143").Value = -1

# Row 198
$ws.Range("B198").Value = 7835232
$ws.Range("E198").Value = "Chippa United"
$ws.Range("F198").Value = "Moroka Swallows"
$ws.Range("G198").Value = 1
$ws.Range("H198").Value = 2
$ws.Range("I198").Value = 1
$ws.Range("J198").Value = 0
$ws.Range("K198").Value = "A"
$ws.Range("L198").Value = 2
$ws.Range("M198").Value = 3.1
$ws.Range("N198").Value = 3.6
$ws.Range("O198").Value = 2.2
$ws.Range("P198").Value = 2.9
$ws.Range("Q198").Value = 3.4
$ws.Range("R198").Value = -0.25
$ws.Range("S198").Value = 1.925
$ws.Range("T198").Value = 1.875
$ws.Range("U198").Value = 1.75
$ws.Range("V198").Value = 1.775
$ws.Range("W198").Value = 2.025
$ws.Range("X198").Value = -1
$ws.Range("Y198").Value = -1
$ws.Range("Z198").Value = 2.4
$ws.Range("AA198").Value = -1
$ws.Range("AB198").Value = 0.875
$ws.Range("AC198").Value = 0.7749999999999999
$ws.Range("AD198").Value = -1

# Row 199
$ws.Range("B199").Value = 7835233
$ws.Range("E199").Value = "Cape Town City"
$ws.Range("F199").Value = "Orlando Pirates"
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 2
$ws.Range("I199").Value = 0
$ws.Range("J199").Value = 0
$ws.Range("K199").Value = "A"
$ws.Range("L199").Value = 3.4
$ws.Range("M199").Value = 3
$ws.Range("N199").Value = 2.15
$ws.Range("O199").Value = 3.1
$ws.Range("P199").Value = 3.2
$ws.Range("Q199").Value = 2.2
$ws.Range("R199").Value = 0.25
$ws.Range("S199").Value = 1.825
$ws.Range("T199").Value = 1.975
$ws.Range("U199").Value = 2.25
$ws.Range("V199").Value = 1.875
$ws.Range("W199").Value = 1.925
$ws.Range("X199").Value = -1
$ws.Range("Y199").Value = -1
$ws.Range("Z199").Value = 1.2
$ws.Range("AA199").Value = -1
$ws.Range("AB199").Value = 0.9750000000000001
$ws.Range("AC199").Value = -0.5
$ws.Range("AD199").Value = 0.4625

# Row 208
$ws.Range("B208").Value = 8163505
$ws.Range("E208").Value = "Orlando Pirates"
$ws.Range("F208").Value = "Chippa United"
$ws.Range("G208").Value = 2
$ws.Range("H208").Value = 0
$ws.Range("I208").Value = 2
$ws.Range("J208").Value = 0
$ws.Range("K208").Value = "H"
$ws.Range("L208").Value = 1.363
$ws.Range("M208").Value = 4.1
$ws.Range("N208").Value = 8.5
$ws.Range("O208").Value = 1.363
$ws.Range("P208").Value = 4.333
$ws.Range("Q208").Value = 8
$ws.Range("R208").Value = -1.25
$ws.Range("S208").Value = 1.8
$ws.Range("T208").Value = 2
$ws.Range("U208").Value = 2.75
$ws.Range("V208").Value = 2
$ws.Range("W208").Value = 1.8
$ws.Range("X208").Value = 0.363
$ws.Range("Y208").Value = -1
$ws.Range("Z208").Value = -1
$ws.Range("AA208").Value = 0.8
$ws.Range("AB208").Value = -1
$ws.Range("AC208").Value = -1
$ws.Range("AD208").Value = 0.8

# Row 211
$ws.Range("B211").Value = 8163498
$ws.Range("E211").Value = "Stellenbosch FC"
$ws.Range("F211").Value = "Amazulu"
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
$ws.Range("I211").Value = 0
$ws.Range("J211").Value = 0
$ws.Range("K211").Value = "D"
$ws.Range("L211").Value = 1.666
$ws.Range("M211").Value = 3.4
$ws.Range("N211").Value = 5
$ws.Range("O211").Value = 1.666
$ws.Range("P211").Value = 3.4
$ws.Range("Q211").Value = 5
$ws.Range("R211").Value = -0.75
$ws.Range("S211").Value = 1.875
$ws.Range("T211").Value = 1.925
$ws.Range("U211").Value = 2
$ws.Range("V211").Value = 1.75
$ws.Range("W211").Value = 2.05
$ws.Range("X211").Value = -1
$ws.Range("Y211").Value = 2.4
$ws.Range("Z211").Value = -1
$ws.Range("AA211").Value = -1
$ws.Range("AB211").Value = 0.925
$ws.Range("AC211").Value = -1
$ws.Range("AD211").Value = 1.05

# Row 226
$ws.Range("B226").Value = 8217791
$ws.Range("E226").Value = "TS Galaxy"
$ws.Range("F226").Value = "Orlando Pirates"
$ws.Range("G226").Value = 1
$ws.Range("H226").Value = 0
$ws.Range("I226").Value = 0
$ws.Range("J226").Value = 0
$ws.Range("K226").Value = "H"
$ws.Range("L226").Value = 3
$ws.Range("M226").Value = 3.2
$ws.Range("N226").Value = 2.25
$ws.Range("O226").Value = 4.1
$ws.Range("P226").Value = 3.2
$ws.Range("Q226").Value = 1.8
$ws.Range("R226").Value = 0.5
$ws.Range("S226").Value = 1.95
$ws.Range("T226").Value = 1.85
$ws.Range("U226").Value = 2
$ws.Range("V226").Value = 1.825
$ws.Range("W226").Value = 1.975
$ws.Range("X226").Value = 3.1
$ws.Range("Y226").Value = -1
$ws.Range("Z226").Value = -1
$ws.Range("AA226").Value = 0.95
$ws.Range("AB226").Value = -1
$ws.Range("AC226").Value = -1
$ws.Range("AD226").Value = 0.9750000000000001

# Row 227
$ws.Range("B227").Value = 7835246
$ws.Range("E227").Value = "Royal AM FC"
$ws.Range("F227").Value = "Chippa United"
$ws.Range("G227").Value = 0
$ws.Range("H227").Value = 0
$ws.Range("I227").Value = 0
$ws.Range("J227").Value = 0
$ws.Range("K227").Value = "D"
$ws.Range("L227").Value = 2.25
$ws.Range("M227").Value = 3
$ws.Range("N227").Value = 3.2
$ws.Range("O227").Value = 2.3
$ws.Range("P227").Value = 2.9
$ws.Range("Q227").Value = 3.2
$ws.Range("R227").Value = -0.25
$ws.Range("S227").Value = 1.975
$ws.Range("T227").Value = 1.825
$ws.Range("U227").Value = 2
$ws.Range("V227").Value = 1.9
$ws.Range("W227").Value = 1.9
$ws.Range("X227").Value = -1
$ws.Range("Y227").Value = 1.9
$ws.Range("Z227").Value = -1
$ws.Range("AA227").Value = -0.5
$ws.Range("AB227").Value = 0.4125
$ws.Range("AC227").Value = -1
$ws.Range("AD227").Value = 0.8999999999999999

# Row 228
$ws.Range("B228").Value = 7835245
$ws.Range("E228").Value = "Kaizer Chiefs"
$ws.Range("F228").Value = "Polokwane City"
$ws.Range("G228").Value = 0
$ws.Range("H228").Value = 0
$ws.Range("I228").Value = 0
$ws.Range("J228").Value = 0
$ws.Range("K228").Value = "D"
$ws.Range("L228").Value = 1.909
$ws.Range("M228").Value = 3.3
$ws.Range("N228").Value = 3.75
$ws.Range("O228").Value = 1.909
$ws.Range("P228").Value = 3
$ws.Range("Q228").Value = 4.333
$ws.Range("R228").Value = -0.5
$ws.Range("S228").Value = 1.95
$ws.Range("T228").Value = 1.85
$ws.Range("U228").Value = 2
$ws.Range("V228").Value = 1.925
$ws.Range("W228").Value = 1.875
$ws.Range("X228").Value = -1
$ws.Range("Y228").Value = 2
$ws.Range("Z228").Value = -1
$ws.Range("AA228").Value = -1
$ws.Range("AB228").Value = 0.8500000000000001
$ws.Range("AC228").Value = -1
$ws.Range("AD228").Value = 0.875

# Row 236
$ws.Range("B236").Value = 7835247
$ws.Range("E236").Value = "Amazulu"
$ws.Range("F236").Value = "Golden Arrows"
$ws.Range("G236").Value = 1
$ws.Range("H236").Value = 3
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = 1
$ws.Range("K236").Value = "A"
$ws.Range("L236").Value = 2
$ws.Range("M236").Value = 3.1
$ws.Range("N236").Value = 3.4
$ws.Range("O236").Value = 1.75
$ws.Range("P236").Value = 3.5
$ws.Range("Q236").Value = 3.8
$ws.Range("R236").Value = -0.5
$ws.Range("S236").Value = 1.8
$ws.Range("T236").Value = 2
$ws.Range("U236").Value = 2.5
$ws.Range("V236").Value = 1.975
$ws.Range("W236").Value = 1.825
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 2.8
$ws.Range("AA236").Value = -1
$ws.Range("AB236").Value = 1
$ws.Range("AC236").Value = 0.9750000000000001
$ws.Range("AD236").Value = -1

# Row 237
$ws.Range("B237").Value = 7834345
$ws.Range("E237").Value = "Chippa United"
$ws.Range("F237").Value = "Sekhukhune United FC"
$ws.Range("G237").Value = 0
$ws.Range("H237").Value = 1
$ws.Range("I237").Value = 0
$ws.Range("J237").Value = 0
$ws.Range("K237").Value = "A"
$ws.Range("L237").Value = 2.5
$ws.Range("M237").Value = 3
$ws.Range("N237").Value = 2.625
$ws.Range("O237").Value = 2.75
$ws.Range("P237").Value = 2.875
$ws.Range("Q237").Value = 2.45
$ws.Range("R237").Value = 0
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 1.8
$ws.Range("U237").Value = 2
$ws.Range("V237").Value = 1.875
$ws.Range("W237").Value = 1.925
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 1.45
$ws.Range("AA237").Value = -1
$ws.Range("AB237").Value = 0.8
$ws.Range("AC237").Value = -1
$ws.Range("AD237").Value = 0.925
